$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

$tbl.ListRows.Add() | Out-Null
$r = 259
$ws.Range("A$r").Value = "3493-kvido-dvoulist.bnl"
$ws.Range("B$r").Value = "cz"
$ws.Range("C$r").Value = 14982656
$ws.Range("D$r").Value = "060cc9b5afe16788e2d1fe5f5ba8f12b"
$ws.Range("E$r").Value = "0x00000DA5"
$ws.Range("F$r").Value = "0x00000002"
$ws.Range("G$r").Value = "0x0000"
$ws.Range("H$r").Value = "0x2EE1"
$ws.Range("I$r").Value = "0x002C"
$ws.Range("J$r").Value = "128kbps CBR stereo/44.1kHz: 86%"

$tbl.ListRows.Add() | Out-Null
$r = 260
$ws.Range("A$r").Value = "anglictina-pro-samouky-3328.bnl"
$ws.Range("B$r").Value = "cz"
$ws.Range("C$r").Value = 463397888
$ws.Range("D$r").Value = "8b89c922c0c1588c3710cfca45064413"
$ws.Range("E$r").Value = "0x00000D00"
$ws.Range("F$r").Value = "0x00000009"
$ws.Range("G$r").Value = "0x0000"
$ws.Range("H$r").Value = "0x4255"
$ws.Range("I$r").Value = "0x22EF"
$ws.Range("J$r").Value = "96kbps CBR mono/44.1kHz: 99%"

$tbl.ListRows.Add() | Out-Null
$r = 261
$ws.Range("A$r").Value = "casopis6-3336.bnl"
$ws.Range("B$r").Value = "cz"
$ws.Range("C$r").Value = 100908544
$ws.Range("D$r").Value = "643d0232aad369d8c59aa800e507ad69"
$ws.Range("E$r").Value = "0x00000D08"
$ws.Range("F$r").Value = "0x00000006"
$ws.Range("G$r").Value = "0x0000"
$ws.Range("H$r").Value = "0x2F4D"
$ws.Range("I$r").Value = "0x045C"
$ws.Range("J$r").Value = "96kbps CBR mono/44.1kHz: 94%"

Write-Host "Done adding rows"
Write-Host "Table ref: " $tbl.Range.Address()

$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add($ws.Range("E1:E261"))
$tbl.Sort.Header = 1
$tbl.Sort.Apply()

Write-Host "Sorted"

$ws.Range("J262").Select()
